$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4716

$ws.Range("H132").Value = 1958.8182
$ws.Range("I132").Value = 1154.8
$ws.Range("K132").Value = 3464.4
$ws.Range("M132").Value = -934.3999999999996

$ws.Range("H135").Value = 701.625
$ws.Range("I135").Value = 659
$ws.Range("K135").Value = 5931
$ws.Range("M135").Value = -3396

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1246.8948
$ws.Range("I2").Value = 1214.1333
$ws.Range("K2").Value = 1214.1333
$ws.Range("M2").Value = -1101.1333

$ws.Range("H32").Value = 3845.1875
$ws.Range("I32").Value = 3845.1875
$ws.Range("K32").Value = 3845.1875
$ws.Range("M32").Value = -3558.1875

$ws.Range("H63").Value = 13648.571
$ws.Range("I63").Value = 3800
$ws.Range("K63").Value = 3800
$ws.Range("M63").Value = -3114

$ws.Range("H66").Value = 13648.571
$ws.Range("I66").Value = 3800
$ws.Range("K66").Value = 19000
$ws.Range("M66").Value = -15568

$ws.Range("H116").Value = 1246.8948
$ws.Range("I116").Value = 1214.1333
$ws.Range("K116").Value = 1214.1333
$ws.Range("M116").Value = 1079.8667

$ws.Range("H132").Value = 3766.5557
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057

$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -155060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1246.8948
$ws.Range("I3").Value = 1214.1333
$ws.Range("K3").Value = 1214.1333
$ws.Range("M3").Value = -1100.1333

$ws.Range("H80").Value = 638.7273
$ws.Range("I80").Value = 101
$ws.Range("J80").Value = 946
$ws.Range("K80").Value = 101
$ws.Range("L80").Value = 946
$ws.Range("M80").Value = 897
$ws.Range("N80").Value = -2942

$ws.Range("H82").Value = 24997.75
$ws.Range("I82").Value = 15997.833
$ws.Range("J82").Value = 51997.5
$ws.Range("K82").Value = 15997.833
$ws.Range("L82").Value = 51997.5
$ws.Range("M82").Value = -15614.833
$ws.Range("N82").Value = -52763.5

$ws.Range("H83").Value = 638.7273
$ws.Range("I83").Value = 101
$ws.Range("J83").Value = 946
$ws.Range("K83").Value = 505
$ws.Range("L83").Value = 4730
$ws.Range("M83").Value = 4487
$ws.Range("N83").Value = -14714

$ws.Range("H85").Value = 24997.75
$ws.Range("I85").Value = 15997.833
$ws.Range("J85").Value = 51997.5
$ws.Range("K85").Value = 15997.833
$ws.Range("L85").Value = 51997.5
$ws.Range("M85").Value = -14671.833
$ws.Range("N85").Value = -54649.5

$ws.Range("H105").Value = 2983.4
$ws.Range("I105").Value = 2009
$ws.Range("J105").Value = 3227
$ws.Range("K105").Value = 2009
$ws.Range("L105").Value = 3227
$ws.Range("M105").Value = -262
$ws.Range("N105").Value = -6721

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -407

$ws.Range("H86").Value = 55000
$ws.Range("I86").Value = 55000
$ws.Range("K86").Value = 55000
$ws.Range("M86").Value = -53877

$ws.Range("H89").Value = 55000
$ws.Range("I89").Value = 55000
$ws.Range("K89").Value = 275000
$ws.Range("M89").Value = -269384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1900
$ws.Range("J4").Value = 2250
$ws.Range("L4").Value = 6750
$ws.Range("N4").Value = -6974

$ws.Range("H74").Value = 7636.364
$ws.Range("J74").Value = 7636.364
$ws.Range("L74").Value = 22909.092
$ws.Range("N74").Value = -25031.092

$ws.Range("H77").Value = 7636.364
$ws.Range("J77").Value = 7636.364
$ws.Range("L77").Value = 68727.276
$ws.Range("N77").Value = -79335.276

$ws.Range("H81").Value = 4041.3333
$ws.Range("I81").Value = 2624.6667
$ws.Range("J81").Value = 4749.6665
$ws.Range("K81").Value = 7874.000100000001
$ws.Range("L81").Value = 14248.9995
$ws.Range("M81").Value = -6751.000100000001
$ws.Range("N81").Value = -16494.9995

$ws.Range("H84").Value = 4041.3333
$ws.Range("I84").Value = 2624.6667
$ws.Range("J84").Value = 4749.6665
$ws.Range("K84").Value = 23622.0003
$ws.Range("L84").Value = 42746.9985
$ws.Range("M84").Value = -18006.0003
$ws.Range("N84").Value = -53978.9985

$ws.Range("H132").Value = 3123.5
$ws.Range("J132").Value = 2994
$ws.Range("L132").Value = 26946
$ws.Range("N132").Value = -32006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -31192

$ws.Range("H46").Value = 39945.5
$ws.Range("I46").Value = 39945
$ws.Range("K46").Value = 39945
$ws.Range("M46").Value = -39789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1354.4286
$ws.Range("I16").Value = 1506.2
$ws.Range("J16").Value = 975
$ws.Range("K16").Value = 1506.2
$ws.Range("L16").Value = 975
$ws.Range("M16").Value = -1336.2
$ws.Range("N16").Value = -1315

$ws.Range("H42").Value = 130014
$ws.Range("J42").Value = 130014
$ws.Range("L42").Value = 130014
$ws.Range("N42").Value = -131140

$ws.Range("H49").Value = 130014
$ws.Range("J49").Value = 130014
$ws.Range("L49").Value = 130014
$ws.Range("N49").Value = -130308

$ws.Range("H68").Value = 1922
$ws.Range("I68").Value = 1922
$ws.Range("K68").Value = 1922
$ws.Range("M68").Value = -1173

$ws.Range("H71").Value = 1922
$ws.Range("I71").Value = 1922
$ws.Range("K71").Value = 9610
$ws.Range("M71").Value = -5866

$ws.Range("H132").Value = 3096.353
$ws.Range("I132").Value = 2280
$ws.Range("K132").Value = 6840
$ws.Range("M132").Value = -4310

$ws.Range("H136").Value = 2210.7
$ws.Range("I136").Value = 2011.8889
$ws.Range("K136").Value = 6035.6667
$ws.Range("M136").Value = -3485.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 89999
$ws.Range("J70").Value = 89999
$ws.Range("L70").Value = 89999
$ws.Range("N70").Value = -90629

$ws.Range("H73").Value = 89999
$ws.Range("J73").Value = 89999
$ws.Range("L73").Value = 89999
$ws.Range("N73").Value = -92183

$ws.Range("H132").Value = 2727.6667
$ws.Range("I132").Value = 1970.6666
$ws.Range("J132").Value = 4998.6665
$ws.Range("K132").Value = 5911.9998
$ws.Range("L132").Value = 14995.9995
$ws.Range("M132").Value = -3381.9998
$ws.Range("N132").Value = -20055.9995

$ws.Range("H136").Value = 4002.5
$ws.Range("I136").Value = 4125.3
$ws.Range("K136").Value = 12375.9
$ws.Range("M136").Value = -9825.900000000001
